$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @(
    "Andar de um lado para o outro ",
    "Comportamentos como se procurasse alguma coisa ",
    "Comportamentos de examinar atentamente ",
    "Entrar em local proibido ",
    "Hiperatividade ",
    "Incapacidade de localizar marcos significativos em um ambiente familiar ",
    "Locomoção ao acaso ",
    "Locomoção inquieta ",
    "Locomoção persistente em busca de alguma coisa ",
    "Locomoção por espaços não autorizados ou privados ",
    "Locomoção que não pode ser facilmente dissuadida ",
    "Locomoção resultando em abandono não intencional de um local ",
    "Longos períodos de locomoção sem destino aparente ",
    "Movimento contínuo de um lugar a outro ",
    "Movimento frequente de um lugar a outro ",
    "Perder-se ",
    "Períodos de locomoção intercalados com períodos de não locomoção (p. ex., sentar, levantar-se, dormir) ",
    "Repetir movimentos do cuidador "
)

$startRow = 322
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "###"
    $ws.Cells.Item($r, 2).Value = 58
    $ws.Cells.Item($r, 3).Value = "$$$"
    $ws.Cells.Item($r, 4).Value = $newValues[$i]
    $ws.Cells.Item($r, 5).Value = "%%%"
}

$null = $ws.Range("D321").Select()
